$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update MAPE, MAE, RSME values in row 2 to the newly computed results
$ws.Range("D2").Value = 0.1394304711723751
$ws.Range("E2").Value = 9.759360307772768
$ws.Range("F2").Value = 26.43340715095914
